$wb = $excel.ActiveWorkbook

# --- Rename Sheet1 -> Sources, add new WordFreqList sheet after it ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sources"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "WordFreqList"

# --- Headers ---
$ws2.Range("A1").Value = "TSCorpus"
$ws2.Range("B1").Value = "Word"
$ws2.Range("C1").Value = "Freq"
$ws2.Range("E1").Value = "my-mini-corpus"
$ws2.Range("F1").Value = "Word"
$ws2.Range("G1").Value = "Freq"

# --- TSCorpus word-frequency list (columns B:C) ---
$words = @("ve","bir","bu","da","de","için","ile","çok","olarak","daha")
$freqs = @(8910007,8185200,5055490,2943937,2776962,2524628,1936336,1870108,1520288,1440179)

# --- my-mini-corpus word-frequency list (columns F:G) ---
$miniWords = @("ve","bir","bu","da","de","ne","boğaziçi","olarak","için","bir")
$miniFreqs = @(276,228,166,87,72,55,54,54,49,46)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws2.Range("B$row").Value = $words[$i]
    $ws2.Range("C$row").Value = $freqs[$i]
    $ws2.Range("F$row").Value = $miniWords[$i]
    $ws2.Range("G$row").Value = $miniFreqs[$i]
}

# --- Build the rotated/centered corpus-name label style once on a scratch
#     cell, then copy its formats onto the two label columns so only a
#     single new cell style gets interned (rather than one per property). ---
$tmpl = $ws2.Range("Z1")
$tmpl.HorizontalAlignment = -4108
$tmpl.VerticalAlignment = -4108
$tmpl.Orientation = 90

$tmpl.Copy()
$ws2.Range("A1:A11").PasteSpecial(-4122)
$ws2.Range("E1:E11").PasteSpecial(-4122)
$tmpl.Clear()

$ws2.Range("A1:A11").Merge()
$ws2.Range("E1:E11").Merge()

$ws2.Range("A1").Select()

Write-Host "done"
